$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.060397028923035
$ws.Range("B1").Value = 1.928446292877197
$ws.Range("C1").Value = 3.117414951324463
$ws.Range("D1").Value = 2.107598066329956
$ws.Range("E1").Value = 0.6939291954040527
